# Updated cryptos list (GitHub Actions scheduled refresh).
# Price (col D) and Volume(1h) (col E) values are plain text cells in the
# source sheet, so number-looking prices are written with a leading
# apostrophe to force text storage, then the style is reset to "Normal" so
# no stray quote-prefix formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.683.39'
$ws.Range('E2').Value = '  +6.96%  '
$ws.Range('D3').Value = '2.398.46'
$ws.Range('E3').Value = '  +3.98%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'113.96"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.30%  '
$ws.Range('D6').Value = "'318.82"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('E7').Value = '  +2.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +3.21%  '
$ws.Range('D10').Value = "'42.07"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.67%  '
$ws.Range('D11').Value = "'0.0929"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.47%  '
$ws.Range('D12').Value = "'8.73"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.16%  '
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('E15').Value = '  +3.54%  '
$ws.Range('D16').Value = '2.761.95'
$ws.Range('E16').Value = '  +3.97%  '
$ws.Range('D17').Value = '2.404.10'
$ws.Range('E17').Value = '  +4.44%  '
$ws.Range('D18').Value = '45.618.49'
$ws.Range('E18').Value = '  +6.42%  '
$ws.Range('D19').Value = "'7.49"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').Value = "'74.75"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = "'3.55"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.56%  '
$ws.Range('D24').Value = "'264.29"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('E25').Value = '  +5.11%  '
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('D27').Value = "'7.63"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.31%  '
$ws.Range('D28').Value = "'11.34"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('D30').Value = "'39.31"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.30%  '
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('D32').Value = "'0.0968"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +12.08%  '
$ws.Range('D33').Value = "'172.58"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.49%  '
$ws.Range('E34').Value = '  +4.87%  '
$ws.Range('E36').Value = '  +6.63%  '
$ws.Range('D37').Value = "'0.118"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('D38').Value = "'4.13"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.00%  '
$ws.Range('E39').Value = '  +8.02%  '
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('E41').Value = '  +12.02%  '
$ws.Range('B42').Value = 'BitcoinSV'
$ws.Range('C42').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D42').Value = "'101.12"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.82%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = "'0.242"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.57%  '
$ws.Range('D44').Value = "'13.58"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.59%  '
$ws.Range('D45').Value = "'72.26"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').Value = "'87.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +14.61%  '
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = "'115.63"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.63%  '
$ws.Range('E49').Value = '  +10.71%  '
$ws.Range('D50').Value = "'9.48"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.68%  '
$ws.Range('D51').Value = '1.669.28'
$ws.Range('E51').Value = '  -3.05%  '
